$wb = $excel.ActiveWorkbook

# --- Translation sheet: add menu + screen transition text ids (rows 4-6) ---
$wsTrans = $wb.Worksheets.Item("Translation")

$wsTrans.Range("B4").Value = "SingleUseId2"
$wsTrans.Range("C4").Value = "Default"
$wsTrans.Range("D4").Value = "Center"
$wsTrans.Range("E4").Value = "LTR"
$wsTrans.Range("F4").Value = "Wyjdź"

$wsTrans.Range("B5").Value = "SingleUseId1"
$wsTrans.Range("C5").Value = "Default"
$wsTrans.Range("D5").Value = "Center"
$wsTrans.Range("E5").Value = "LTR"
$wsTrans.Range("F5").Value = "Nowa Gra"

$wsTrans.Range("B6").Value = "SingleUseId3"
$wsTrans.Range("C6").Value = "Typography_00"
$wsTrans.Range("D6").Value = "Left"
$wsTrans.Range("E6").Value = "LTR"
$wsTrans.Range("F6").Value = "FlaBi"

# --- Typography sheet: add a new font entry (row 7) used by the new texts ---
$wsTypo = $wb.Worksheets.Item("Typography")

$wsTypo.Range("B7:F7").Style = "Normal"
$wsTypo.Range("B7").Value = "Typography_00"
$wsTypo.Range("C7").Value = "comicbd.ttf"
$wsTypo.Range("D7").Value = 83
$wsTypo.Range("E7").Value = 4
$wsTypo.Range("F7").Value = "?"

# touch G7:J7 so they materialize as empty cells, matching the other data rows
$wsTypo.Range("G7").Font.Bold = $false
$wsTypo.Range("H7").Font.Bold = $false
$wsTypo.Range("I7").Font.Bold = $false
$wsTypo.Range("J7").Font.Bold = $false
